# Incorporate updated data from upstream processes through 2024
#
# The chart "Distributed Energy Resources Capacity by Year Opened - London"
# plots the Solar column (column E) against Open year (column A) on
# Sheet1. The 2022 (row 24) and 2024 (row 26) Solar capacity figures are
# revised upward to reflect updated upstream data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2022 Solar capacity: 37.72 -> 40.33
$ws.Range("E24").Value = 40.33

# 2024 Solar capacity: 45.72 -> 79.56
$ws.Range("E26").Value = 79.56
